# Acceptance Test Plan.xlsx update
# - Fills in the Sprint 3 "Pass/Fail" (col G) and "Tester initials; date; comments" (col H)
#   columns for all rows that already have Sprint 1/2 test data (rows 2-30 and 32),
#   recording a fresh round of testing on 11/24 that all passed.
# - Replaces the old "As a Player and Spectator ... watching the game" story (row 31)
#   with a new spectator story.
# - Appends two new user stories (rows 33-34) with their own Sprint 3 results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

$passLabel = "Pass"
$comment = "JL; 11/24"

# Rows 2-30 and 32 already contained a user story / acceptance criterion row;
# row 31 is being replaced and rows 33-34 are brand new, handled separately below.
$existingRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,32)
foreach ($r in $existingRows) {
    $ws.Cells.Item($r, 7).Value = $passLabel   # column G
    $ws.Cells.Item($r, 8).Value = $comment     # column H
}

# Row 31: story text changes from the old "Player and Spectator" viewer-count story
# to a new "spectator" story, and also gets Sprint 3 results like the rest.
$ws.Cells.Item(31, 1).Value = "As a spectator, I do not want to be asked to play in a game"
$ws.Cells.Item(31, 7).Value = $passLabel
$ws.Cells.Item(31, 8).Value = $comment

# New rows 33 and 34: brand-new user stories with Sprint 3 results filled in.
$ws.Cells.Item(33, 1).Value = "As a player, I want to play with an AI"
$ws.Cells.Item(33, 7).Value = $passLabel
$ws.Cells.Item(33, 8).Value = $comment

$ws.Cells.Item(34, 1).Value = "As a player, I want to see a message of how I won or lost when I lose or win a game"
$ws.Cells.Item(34, 7).Value = $passLabel
$ws.Cells.Item(34, 8).Value = $comment

# Update the saved cursor/selection to reflect where editing left off.
$ws.Activate()
$ws.Range("H35").Select()

Write-Output "done"
